# edit.ps1 - apply cryptos.xlsx price/volume updates described by the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value reads as a plain decimal number (e.g. "247.15").
# Excel would silently convert these to numeric values and drop formatting
# (e.g. trailing zeros) unless we force a Text number format first, matching
# the original inlineStr/text storage of these price cells.
$textCells = @("D5", "D6", "D7", "D10", "D11", "D12", "D13", "D15", "D16", "D21", "D23", "D24", "D28", "D30", "D31", "D32", "D33", "D36", "D37", "D39", "D42", "D44", "D45", "D49", "D50", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply the updated cell values (coin prices and 1h volume deltas), and the
# TrustWalletToken / FTXToken row swap (rows 49-50).
$ws.Range('D2').Value = '42.020.28'
$ws.Range('E2').Value = '  -1.65%  '
$ws.Range('D3').Value = '2.243.66'
$ws.Range('E3').Value = '  -2.06%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '247.15'
$ws.Range('E5').Value = '  -1.87%  '
$ws.Range('D6').Value = '0.631'
$ws.Range('E6').Value = '  -1.06%  '
$ws.Range('D7').Value = '75.26'
$ws.Range('E7').Value = '  +1.60%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  -2.43%  '
$ws.Range('D10').Value = '40.42'
$ws.Range('E10').Value = '  +2.66%  '
$ws.Range('D11').Value = '0.0953'
$ws.Range('E11').Value = '  -3.71%  '
$ws.Range('D12').Value = '7.19'
$ws.Range('E12').Value = '  -2.02%  '
$ws.Range('D13').Value = '0.103'
$ws.Range('E13').Value = '  -2.52%  '
$ws.Range('D14').Value = '2.579.57'
$ws.Range('E14').Value = '  -2.06%  '
$ws.Range('D15').Value = '14.88'
$ws.Range('E15').Value = '  -3.51%  '
$ws.Range('D16').Value = '0.864'
$ws.Range('E16').Value = '  -1.78%  '
$ws.Range('D17').Value = '2.257.54'
$ws.Range('E17').Value = '  -1.33%  '
$ws.Range('D18').Value = '41.947.07'
$ws.Range('E18').Value = '  -1.68%  '
$ws.Range('D19').Value = '0.0₃0981'
$ws.Range('E19').Value = '  -1.69%  '
$ws.Range('E20').Value = '  -2.44%  '
$ws.Range('D21').Value = '71.70'
$ws.Range('E21').Value = '  -1.14%  '
$ws.Range('E22').Value = '  +0.51%  '
$ws.Range('D23').Value = '231.29'
$ws.Range('E23').Value = '  -1.77%  '
$ws.Range('D24').Value = '11.45'
$ws.Range('E24').Value = '  +0.11%  '
$ws.Range('E25').Value = '  +0.09%  '
$ws.Range('E26').Value = '  -5.77%  '
$ws.Range('E27').Value = '  -4.59%  '
$ws.Range('D28').Value = '7.14'
$ws.Range('E28').Value = '  +11.31%  '
$ws.Range('E29').Value = '  -1.41%  '
$ws.Range('D30').Value = '168.89'
$ws.Range('E30').Value = '  +1.15%  '
$ws.Range('D31').Value = '20.57'
$ws.Range('E31').Value = '  -2.42%  '
$ws.Range('D32').Value = '33.74'
$ws.Range('E32').Value = '  +4.60%  '
$ws.Range('D33').Value = '0.0844'
$ws.Range('E33').Value = '  +2.94%  '
$ws.Range('E34').Value = '  -4.74%  '
$ws.Range('E35').Value = '  -0.15%  '
$ws.Range('D36').Value = '4.55'
$ws.Range('E36').Value = '  -3.81%  '
$ws.Range('D37').Value = '4.87'
$ws.Range('E37').Value = '  +2.37%  '
$ws.Range('E38').Value = '  -2.42%  '
$ws.Range('D39').Value = '13.48'
$ws.Range('E39').Value = '  -6.84%  '
$ws.Range('E40').Value = '  -0.42%  '
$ws.Range('E41').Value = '  -6.65%  '
$ws.Range('D42').Value = '112.10'
$ws.Range('E42').Value = '  +13.97%  '
$ws.Range('E43').Value = '  -4.36%  '
$ws.Range('D44').Value = '60.55'
$ws.Range('E44').Value = '  -2.32%  '
$ws.Range('D45').Value = '8.77'
$ws.Range('E45').Value = '  -3.97%  '
$ws.Range('E46').Value = '  -2.50%  '
$ws.Range('E47').Value = '  -0.37%  '
$ws.Range('E48').Value = '  -3.79%  '
$ws.Range('B49').Value = 'FTXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D49').Value = '4.33'
$ws.Range('E49').Value = '  -10.76%  '
$ws.Range('B50').Value = 'TrustWalletToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D50').Value = '1.17'
$ws.Range('E50').Value = '  -1.93%  '
$ws.Range('D51').Value = '4.20'
$ws.Range('E51').Value = '  -1.66%  '
